# [ADD] New version of BrakeAbaqus
# Updates the brake-sizing parameters on "Plan1":
#   - l_p  (B3): 2.67 -> 3.32
#   - Amc  (B4/C4): 1.98 -> 3.87
# and touches a new cell (J11) so the sheet's used range grows out to
# column J, matching the extra (currently empty but formatted) cell that
# shows up in the new version of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated parameter values ------------------------------------------------
$ws.Range("B3").Value = 3.32
$ws.Range("B4").Value = 3.87
$ws.Range("C4").Value = 3.87

# --- new, formatted-but-empty cell at J11 ------------------------------------
# Round-tripping the font name forces Excel to record an explicit (if
# no-op) font attribute on the cell, which is what materializes J11 with
# its own style and expands the sheet dimension/row spans out to column J.
$ws.Range("J11").Font.Name = "Arial"
$ws.Range("J11").Font.Name = "Calibri"

# --- selection now sits on the newly touched cell ----------------------------
$ws.Range("J11").Select()
